# Update legacy GSC export data: the daily rolling window advanced by one
# day. Every date in column A shifts forward one day (2025-10-14 drops off
# the front, 2026-01-12 is appended at the back), and column C's per-day
# values shift up to match (each row now shows the following day's figure),
# with the newly-appended last day getting a fresh value of 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 91

# Capture the current (pre-edit) column A/C values before we start overwriting.
$oldA = @{}
$oldC = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldA[$r] = $ws.Cells.Item($r, 1).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
}

# The dates are stored as plain text (yyyy-MM-dd), not real Excel date
# serials. Force the column to a text number format before writing so
# Excel's auto-conversion doesn't turn the new date-looking strings into
# date serials, then clear that temporary formatting again so the cells
# end up with no explicit style, same as the source file.
$colA = $ws.Range("A" + $firstRow + ":A" + $lastRow)
$colA.NumberFormat = "@"

# Column A: shift every date forward by one day (row r becomes what row
# r+1 used to be), and the final row gets the day after the old last date.
for ($r = $firstRow; $r -lt $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $oldA[$r + 1]
}
$lastDate = [DateTime]::ParseExact($oldA[$lastRow], "yyyy-MM-dd", $null)
$newDate = $lastDate.AddDays(1)
$ws.Cells.Item($lastRow, 1).Value = $newDate.ToString("yyyy-MM-dd")

$colA.ClearFormats()

# Column C: shift each day's count up by one row to match the new dates;
# the newly appended last day has no prior figure, so it starts at 0.
for ($r = $firstRow; $r -lt $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $oldC[$r + 1]
}
$ws.Cells.Item($lastRow, 3).Value = 0
